$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.726
$ws.Range("B3").Value = 6.427
$ws.Range("C5").Value = -12.836
$ws.Range("E5").Value = 13.034
$ws.Range("E9").Value = 13.649
$ws.Range("E11").Value = 13.071
$ws.Range("B14").Value = 7.166000000000001
$ws.Range("B16").Value = 5.952999999999999
$ws.Range("C16").Value = -12.093
$ws.Range("E17").Value = 13.653
$ws.Range("B21").Value = 6.6
$ws.Range("E21").Value = 12.694
$ws.Range("B23").Value = 6.609999999999999
$ws.Range("B25").Value = 6.355
